# Updated cryptos list on Mon Jun 24 14:52:04 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto tracker sheet, and fixes two pairs of rows whose ranking order
# swapped (WrappedBTC/Avalanche, USDe/Kaspa, SuiNetwork/InjectiveProtocol).
#
# Note: several price values look like plain numbers (e.g. "7.31", "1.00").
# Excel would normally auto-convert such strings to numeric cells, which
# would both change the cell type and lose formatting (e.g. "1.00" -> 1).
# To preserve them as text exactly like the source data, those values are
# assigned with a leading apostrophe, Excel's standard "treat as text"
# quote-prefix - the apostrophe itself is not stored in the cell value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.168.86'
$ws.Range("E2").Value = '  -4.75%  '
$ws.Range("D3").Value = '3.302.58'
$ws.Range("E3").Value = '  -5.33%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''565.56'
$ws.Range("E5").Value = '  -3.75%  '
$ws.Range("D6").Value = '''128.13'
$ws.Range("E6").Value = '  -3.29%  '
$ws.Range("D8").Value = '3.302.53'
$ws.Range("E8").Value = '  -5.30%  '
$ws.Range("D10").Value = '''7.31'
$ws.Range("E10").Value = '  -5.25%  '
$ws.Range("D11").Value = '''0.118'
$ws.Range("E11").Value = '  -4.23%  '
$ws.Range("D12").Value = '''0.375'
$ws.Range("E12").Value = '  -2.98%  '
$ws.Range("D13").Value = '3.875.74'
$ws.Range("E13").Value = '  -5.25%  '
$ws.Range("D14").Value = '''0.119'
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").Value = '3.311.47'
$ws.Range("E15").Value = '  -5.24%  '
$ws.Range("D16").Value = '''0.0000168'
$ws.Range("E16").Value = '  -5.35%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '61.291.48'
$ws.Range("E17").Value = '  -4.53%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").Value = '''24.46'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("E19").Value = '  -1.46%  '
$ws.Range("D20").Value = '''13.33'
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("D21").Value = '''8.96'
$ws.Range("E21").Value = '  -10.32%  '
$ws.Range("D22").Value = '''355.14'
$ws.Range("E22").Value = '  -7.67%  '
$ws.Range("D23").Value = '''0.554'
$ws.Range("E23").Value = '  -3.89%  '
$ws.Range("D25").Value = '3.438.59'
$ws.Range("E25").Value = '  -5.33%  '
$ws.Range("D26").Value = '''69.27'
$ws.Range("E26").Value = '  -6.87%  '
$ws.Range("D27").Value = '''0.0000107'
$ws.Range("E27").Value = '  -5.68%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").Value = '''7.20'
$ws.Range("D30").Value = '''1.44'
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("D31").Value = '''7.82'
$ws.Range("E31").Value = '  -1.93%  '
$ws.Range("D32").Value = '''2.10'
$ws.Range("E32").Value = '  -6.01%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '''1.00'
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '''0.149'
$ws.Range("E34").Value = '  -3.05%  '
$ws.Range("D35").Value = '3.332.67'
$ws.Range("E35").Value = '  -5.37%  '
$ws.Range("D36").Value = '''22.54'
$ws.Range("E36").Value = '  -2.18%  '
$ws.Range("D37").Value = '''5.25'
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("D38").Value = '''6.78'
$ws.Range("E38").Value = '  -0.80%  '
$ws.Range("D39").Value = '''161.23'
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("D40").Value = '''1.48'
$ws.Range("E40").Value = '  -3.16%  '
$ws.Range("D41").Value = '''0.0757'
$ws.Range("E41").Value = '  -2.93%  '
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").Value = '''4.37'
$ws.Range("E43").Value = '  +0.33%  '
$ws.Range("D44").Value = '''41.03'
$ws.Range("E44").Value = '  -1.64%  '
$ws.Range("D45").Value = '''0.742'
$ws.Range("E45").Value = '  -7.66%  '
$ws.Range("E46").Value = '  -4.76%  '
$ws.Range("D47").Value = '''1.55'
$ws.Range("E47").Value = '  -4.93%  '
$ws.Range("D48").Value = '''22.28'
$ws.Range("E48").Value = '  -8.11%  '
$ws.Range("D49").Value = '''6.69'
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").Value = '''0.853'
$ws.Range("E50").Value = '  -9.04%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '''21.26'
$ws.Range("E51").Value = '  +2.88%  '
